$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.587.98"
$ws.Range("E2").Value = "  +6.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.559.35"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "413.10"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.62"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.552.73"
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.765"
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.171"
$ws.Range("E11").Value = "  +20.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000310"
$ws.Range("E12").Value = "  +42.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.93"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.78"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.126.26"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.02"
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.588.09"
$ws.Range("E18").Value = "  +3.77%  "
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.538.50"
$ws.Range("E20").Value = "  +6.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.17"
$ws.Range("E21").Value = "  -4.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "443.57"
$ws.Range("E22").Value = "  -4.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.85"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.11"
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.90"
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.30"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").Value = "  -7.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.17"
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.23"
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.24"
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.158"
$ws.Range("E34").Value = "  -5.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.07"
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.23"
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0486"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0702"
$ws.Range("E39").Value = "  +23.68%  "
$ws.Range("E40").Value = "  +9.52%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.94"
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.70"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "146.47"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.26"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("E48").Value = "  -6.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.06"
$ws.Range("E49").Value = "  +7.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("E50").Value = "  -7.24%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.54"
$ws.Range("E51").Value = "  +9.03%  "
